$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows right before row 19. This pushes the existing
# "last row" (VICTOR BELLIDO duplicate, special border styling) down to row 21,
# leaving fresh rows 19 and 20 to fill with new worker data.
$ws.Rows("19:20").Insert()

# The inserted rows come back with generic default formatting; copy the
# "interior" row formatting (row 18) into them so they match the existing
# table look (thin borders, centered, currency format, etc.).
$ws.Range("B18:J18").Copy()
$ws.Range("B19:J19").PasteSpecial(-4122)
$ws.Range("B18:J18").Copy()
$ws.Range("B20:J20").PasteSpecial(-4122)

# --- Row 17 (VICTOR BELLIDO RIVERA): period corrected from 2507 to 2505 ---
$ws.Range("E17").Value = "2505"

# --- Row 18: new worker JORGE NUÑEZ SUAREZ ---
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "9095580"
$ws.Range("D18").Value = "JORGE NUÑEZ SUAREZ"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 63980
$ws.Range("G18").Value = 1599520

# --- Row 19: new worker WILFRIDO AYALA PEREZ ---
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1128060506"
$ws.Range("D19").Value = "WILFRIDO AYALA PEREZ"
$ws.Range("E19").Value = "2508"
$ws.Range("F19").Value = 63980
$ws.Range("G19").Value = 1599520

# --- Row 20: new worker ROBERTO CARLOS GAMARRA LORA ---
$ws.Range("B20").Value = "CC"
$ws.Range("C20").Value = "1005627516"
$ws.Range("D20").Value = "ROBERTO CARLOS GAMARRA LORA"
$ws.Range("E20").Value = "2508"
$ws.Range("F20").Value = 83200
$ws.Range("G20").Value = 2080000

# --- Row 21 (previously row 19): new worker JOSE DEMETRIO CARRASQUILLA MANJARRES ---
$ws.Range("B21").Value = "CC"
$ws.Range("C21").Value = "1049935981"
$ws.Range("D21").Value = "JOSE DEMETRIO CARRASQUILLA MANJARRES"
$ws.Range("E21").Value = "2508"
$ws.Range("F21").Value = 60000
$ws.Range("G21").Value = 1500000

# --- Summary header updates ---
# Valor Mora total: sum of F16:F21
$ws.Range("E11").Value = 385040
# Cant. Trabajadores now 6, Cant. Periodos now 2 (2505 & 2508)
$ws.Range("C13").Value = 6
$ws.Range("F13").Value = 2

# Column D widens to fit the longest new worker name
$ws.Columns("D").ColumnWidth = 41.90625
